$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 9481
$ws.Range("E2").Value = 27
$ws.Range("F2").Value = 26
$ws.Range("G2").Value = 13
$ws.Range("H2").Value = 7
$ws.Range("I2").Value = 7
$ws.Range("K2").Value = 15904
$ws.Range("L2").Value = 8211
$ws.Range("M2").Value = 7694
$ws.Range("N2").Value = 7694
$ws.Range("P2").Value = 3941
$ws.Range("Q2").Value = -479
$ws.Range("R2").Value = -102
$ws.Range("S2").Value = 707
$ws.Range("T2").Value = 90
$ws.Range("U2").Value = -569
$ws.Range("V2").Value = 5278
$ws.Range("W2").Value = 0.28
$ws.Range("X2").Value = 0.08
$ws.Range("Y2").Value = 0.09
$ws.Range("Z2").Value = 0.05
$ws.Range("AA2").Value = 106.72
$ws.Range("AB2").Value = 100.43
$ws.Range("AC2").Value = 9
$ws.Range("AD2").Value = 454.02
$ws.Range("AE2").Value = 10326
$ws.Range("AF2").Value = 0.4
$ws.Range("AG2").Value = 10
$ws.Range("AH2").Value = 0.24
$ws.Range("AI2").Value = 61.53
$ws.Range("AJ2").Value = 77636164
$ws.Range("J2").ClearContents()
$ws.Range("O2").ClearContents()

# Row 3
$ws.Range("D3").Value = 7310
$ws.Range("E3").Value = -728
$ws.Range("F3").Value = -684
$ws.Range("G3").Value = -967
$ws.Range("H3").Value = -758
$ws.Range("I3").Value = -758
$ws.Range("K3").Value = 15211
$ws.Range("L3").Value = 8312
$ws.Range("M3").Value = 6899
$ws.Range("N3").Value = 6899
$ws.Range("P3").Value = 3941
$ws.Range("Q3").Value = 109
$ws.Range("R3").Value = -680
$ws.Range("S3").Value = 426
$ws.Range("T3").Value = 130
$ws.Range("U3").Value = -22
$ws.Range("V3").Value = 5946
$ws.Range("W3").Value = -9.960000000000001
$ws.Range("X3").Value = -10.36
$ws.Range("Y3").Value = -10.38
$ws.Range("Z3").Value = -4.87
$ws.Range("AA3").Value = 120.48
$ws.Range("AB3").Value = 81.17
$ws.Range("AC3").Value = -961
$ws.Range("AD3").Value = -3.79
$ws.Range("AE3").Value = 9387
$ws.Range("AF3").Value = 0.39
$ws.Range("AG3").Value = 0
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = 78798750
$ws.Range("J3").ClearContents()
$ws.Range("O3").ClearContents()

# Row 4
$ws.Range("D4").Value = 4870
$ws.Range("E4").Value = -582
$ws.Range("F4").Value = -582
$ws.Range("G4").Value = -1016
$ws.Range("H4").Value = -871
$ws.Range("I4").Value = -871
$ws.Range("K4").Value = 12613
$ws.Range("L4").Value = 6584
$ws.Range("M4").Value = 6029
$ws.Range("N4").Value = 6029
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = 3941
$ws.Range("Q4").Value = -242
$ws.Range("R4").Value = 1528
$ws.Range("S4").Value = -1176
$ws.Range("T4").Value = 96
$ws.Range("U4").Value = -338
$ws.Range("V4").Value = 4956
$ws.Range("W4").Value = -11.96
$ws.Range("X4").Value = -17.89
$ws.Range("Y4").Value = -13.48
$ws.Range("Z4").Value = -6.26
$ws.Range("AA4").Value = 109.2
$ws.Range("AB4").Value = 59.18
$ws.Range("AC4").Value = -1106
$ws.Range("AD4").Value = -2.75
$ws.Range("AE4").Value = 8203
$ws.Range("AF4").Value = 0.37
$ws.Range("AG4").Value = 0
$ws.Range("AH4").Value = 0
$ws.Range("AI4").Value = 0
$ws.Range("AJ4").Value = 78803016
$ws.Range("J4").ClearContents()

# Row 5
$ws.Range("D5").Value = 4434
$ws.Range("E5").Value = 28
$ws.Range("F5").Value = 28
$ws.Range("G5").Value = 331
$ws.Range("H5").Value = 206
$ws.Range("I5").Value = 206
$ws.Range("K5").Value = 12735
$ws.Range("L5").Value = 6506
$ws.Range("M5").Value = 6229
$ws.Range("N5").Value = 6229
$ws.Range("O5").Value = 0
$ws.Range("P5").Value = 3941
$ws.Range("Q5").Value = 187
$ws.Range("R5").Value = 357
$ws.Range("S5").Value = -333
$ws.Range("T5").Value = 77
$ws.Range("U5").Value = 110
$ws.Range("V5").Value = 4796
$ws.Range("W5").Value = 0.63
$ws.Range("X5").Value = 4.65
$ws.Range("Y5").Value = 3.36
$ws.Range("Z5").Value = 1.63
$ws.Range("AA5").Value = 104.46
$ws.Range("AB5").Value = 64.16
$ws.Range("AC5").Value = 261
$ws.Range("AD5").Value = 11.08
$ws.Range("AE5").Value = 8474
$ws.Range("AF5").Value = 0.34
$ws.Range("AI5").Value = 0
$ws.Range("AJ5").Value = 78803016
$ws.Range("J5").ClearContents()
$ws.Range("AG5").ClearContents()
$ws.Range("AH5").ClearContents()

# Row 6
$ws.Range("D6").Value = 4459
$ws.Range("E6").Value = 38
$ws.Range("F6").Value = 38
$ws.Range("G6").Value = 191
$ws.Range("H6").Value = 163
$ws.Range("I6").Value = 163
$ws.Range("K6").Value = 12736
$ws.Range("L6").Value = 6355
$ws.Range("M6").Value = 6381
$ws.Range("N6").Value = 6381
$ws.Range("P6").Value = 3941
$ws.Range("Q6").Value = 90
$ws.Range("R6").Value = 45
$ws.Range("S6").Value = -297
$ws.Range("T6").Value = 61
$ws.Range("U6").Value = 29
$ws.Range("V6").Value = 4683
$ws.Range("W6").Value = 0.86
$ws.Range("X6").Value = 3.65
$ws.Range("Y6").Value = 2.59
$ws.Range("Z6").Value = 1.28
$ws.Range("AA6").Value = 99.58
$ws.Range("AB6").Value = 68.05
$ws.Range("AC6").Value = 207
$ws.Range("AD6").Value = 11.99
$ws.Range("AE6").Value = 8683
$ws.Range("AF6").Value = 0.29
$ws.Range("AI6").Value = 0
$ws.Range("AJ6").Value = 78803016
$ws.Range("AG6").ClearContents()
$ws.Range("AH6").ClearContents()

# Row 7
$ws.Range("D7").ClearContents()
$ws.Range("E7").ClearContents()
$ws.Range("G7").ClearContents()
$ws.Range("H7").ClearContents()
$ws.Range("I7").ClearContents()
$ws.Range("K7").ClearContents()
$ws.Range("L7").ClearContents()
$ws.Range("M7").ClearContents()
$ws.Range("N7").ClearContents()
$ws.Range("P7").ClearContents()
$ws.Range("Q7").ClearContents()
$ws.Range("R7").ClearContents()
$ws.Range("S7").ClearContents()
$ws.Range("T7").ClearContents()
$ws.Range("U7").ClearContents()
$ws.Range("W7").ClearContents()
$ws.Range("X7").ClearContents()
$ws.Range("Y7").ClearContents()
$ws.Range("Z7").ClearContents()
$ws.Range("AA7").ClearContents()
$ws.Range("AC7").ClearContents()
$ws.Range("AD7").ClearContents()
$ws.Range("AE7").ClearContents()
$ws.Range("AF7").ClearContents()
$ws.Range("AG7").ClearContents()
$ws.Range("AH7").ClearContents()
$ws.Range("AI7").ClearContents()

# Row 8
$ws.Range("D8").ClearContents()
$ws.Range("E8").ClearContents()
$ws.Range("G8").ClearContents()
$ws.Range("H8").ClearContents()
$ws.Range("I8").ClearContents()
$ws.Range("K8").ClearContents()
$ws.Range("L8").ClearContents()
$ws.Range("M8").ClearContents()
$ws.Range("N8").ClearContents()
$ws.Range("P8").ClearContents()
$ws.Range("Q8").ClearContents()
$ws.Range("R8").ClearContents()
$ws.Range("S8").ClearContents()
$ws.Range("T8").ClearContents()
$ws.Range("U8").ClearContents()
$ws.Range("W8").ClearContents()
$ws.Range("X8").ClearContents()
$ws.Range("Y8").ClearContents()
$ws.Range("Z8").ClearContents()
$ws.Range("AA8").ClearContents()
$ws.Range("AC8").ClearContents()
$ws.Range("AD8").ClearContents()
$ws.Range("AE8").ClearContents()
$ws.Range("AF8").ClearContents()
$ws.Range("AG8").ClearContents()
$ws.Range("AH8").ClearContents()
$ws.Range("AI8").ClearContents()

# Row 9
$ws.Range("D9").ClearContents()
$ws.Range("E9").ClearContents()
$ws.Range("G9").ClearContents()
$ws.Range("H9").ClearContents()
$ws.Range("I9").ClearContents()
$ws.Range("K9").ClearContents()
$ws.Range("L9").ClearContents()
$ws.Range("M9").ClearContents()
$ws.Range("N9").ClearContents()
$ws.Range("P9").ClearContents()
$ws.Range("Q9").ClearContents()
$ws.Range("R9").ClearContents()
$ws.Range("S9").ClearContents()
$ws.Range("T9").ClearContents()
$ws.Range("U9").ClearContents()
$ws.Range("W9").ClearContents()
$ws.Range("X9").ClearContents()
$ws.Range("Y9").ClearContents()
$ws.Range("Z9").ClearContents()
$ws.Range("AA9").ClearContents()
$ws.Range("AC9").ClearContents()
$ws.Range("AD9").ClearContents()
$ws.Range("AE9").ClearContents()
$ws.Range("AF9").ClearContents()
$ws.Range("AG9").ClearContents()
$ws.Range("AH9").ClearContents()
$ws.Range("AI9").ClearContents()
